# Update the "F" column (collection/favorite counts) for a handful of
# events across the four worksheets, matching the regenerated output
# published to gh-pages.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F9").Value = 857
$ws.Range("F11").Value = 401
$ws.Range("F19").Value = 1664
$ws.Range("F31").Value = 2322
$ws.Range("F32").Value = 379
$ws.Range("F37").Value = 179
$ws.Range("F41").Value = 409
$ws.Range("F42").Value = 369

# 演出 (Performances)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 44
$ws.Range("F23").Value = 90

# 本地生活 (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2915
$ws.Range("F6").Value = 298

# 全部类型 (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 44
$ws.Range("F13").Value = 857
$ws.Range("F15").Value = 401
$ws.Range("F22").Value = 298
$ws.Range("F23").Value = 1664
$ws.Range("F38").Value = 2322
$ws.Range("F43").Value = 179
